# Apply the commit's layout / formatting tweaks to the mockup slides.
#
# Shape.Left/Top/Width/Height on this object model are expressed in points
# (1 pt = 12700 EMU) and are stored internally as single-precision floats,
# exactly like real PowerPoint COM automation. To make sure the EMU value
# that ends up in the OOXML round-trips to the exact target (rather than
# being off by one EMU because of float32 truncation), nudge the point
# value up by half an EMU before the division.
function EmuToPt($emu) {
    return ([double]$emu + 0.5) / 12700.0
}

# Font.Color.RGB uses the Windows COLORREF packing (0x00BBGGRR), i.e. the
# low byte is Red, the next is Green, the high byte is Blue - convert from
# a normal "RRGGBB" hex string used in the OOXML srgbClr values.
function RGBColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2),16)
    $g = [Convert]::ToInt32($hex.Substring(2,2),16)
    $b = [Convert]::ToInt32($hex.Substring(4,2),16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 - "Mockup": reposition the screenshot picture (bottom-left one)
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$pic1 = $s1.Shapes.Item(3)
$pic1.Left = EmuToPt 0
$pic1.Top  = EmuToPt 4462650

# ---------------------------------------------------------------------
# Slide 2 - "Movil": nudge the screenshot picture
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$pic2 = $s2.Shapes.Item(2)
$pic2.Left = EmuToPt 1073088
$pic2.Top  = EmuToPt 2844079

# ---------------------------------------------------------------------
# Slide 3 - "Principal": move/resize picture + caption textbox
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$pic3 = $s3.Shapes.Item(2)
$pic3.Left   = EmuToPt 889677
$pic3.Top    = EmuToPt 2312266
$pic3.Width  = EmuToPt 2398468

$txt3 = $s3.Shapes.Item(3)
$txt3.Left   = EmuToPt 3288145
$txt3.Top    = EmuToPt 3279847
$txt3.Height = EmuToPt 954107
$tr3 = $txt3.TextFrame.TextRange
$tr3.ParagraphFormat.Alignment = 2   # ppAlignCenter
$tr3.Font.Size = 28

# ---------------------------------------------------------------------
# Slide 4 - "PC": shrink the picture height slightly
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$pic4 = $s4.Shapes.Item(2)
$pic4.Height = EmuToPt 3870757

# ---------------------------------------------------------------------
# Slide 5 - "Menú": title colour + picture/caption repositioning
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$title5 = $s5.Shapes.Item(1)
$title5.TextFrame.TextRange.Font.Color.RGB = RGBColor "FFC000"

$pic5 = $s5.Shapes.Item(2)
$pic5.Left = EmuToPt 1624177
$pic5.Top  = EmuToPt 2400733

$txt5 = $s5.Shapes.Item(3)
$txt5.Left   = EmuToPt 3934691
$txt5.Top    = EmuToPt 3121892
$txt5.Height = EmuToPt 1815882
$tr5 = $txt5.TextFrame.TextRange
$tr5.ParagraphFormat.Alignment = 2   # ppAlignCenter
$tr5.Font.Size = 28

# ---------------------------------------------------------------------
# Slide 6 - "PC": title colour change
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$title6 = $s6.Shapes.Item(1)
$title6.TextFrame.TextRange.Font.Color.RGB = RGBColor "FFC000"

# ---------------------------------------------------------------------
# Slide 7 - "Noticia": title colour + picture/caption repositioning
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$title7 = $s7.Shapes.Item(1)
$title7.TextFrame.TextRange.Font.Color.RGB = RGBColor "FFC000"

$pic7 = $s7.Shapes.Item(2)
$pic7.Left  = EmuToPt 882193
$pic7.Width = EmuToPt 2165808

$txt7 = $s7.Shapes.Item(3)
$txt7.Left   = EmuToPt 4156363
$txt7.Height = EmuToPt 1815882
$tr7 = $txt7.TextFrame.TextRange
$tr7.ParagraphFormat.Alignment = 2   # ppAlignCenter
$tr7.Font.Size = 28

# ---------------------------------------------------------------------
# Slide 9 - "Creador de Noticias": resize picture, centre caption
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$pic9 = $s9.Shapes.Item(2)
$pic9.Width = EmuToPt 2075329

$txt9 = $s9.Shapes.Item(3)
$tr9 = $txt9.TextFrame.TextRange
$tr9.ParagraphFormat.Alignment = 2   # ppAlignCenter

# ---------------------------------------------------------------------
# Slide 11 - "Comunitario": resize picture + move/resize caption
# ---------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$pic11 = $s11.Shapes.Item(2)
$pic11.Width = EmuToPt 1793437

$txt11 = $s11.Shapes.Item(3)
$txt11.Left   = EmuToPt 4257964
$txt11.Top    = EmuToPt 3113229
$txt11.Height = EmuToPt 1384995
$tr11 = $txt11.TextFrame.TextRange
$tr11.Font.Size = 28
